$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.93"
$ws.Range("E2").Value = "'-0.80%"
$ws.Range("D3").Value = "'35.78"
$ws.Range("E3").Value = "'-0.53%"
$ws.Range("D4").Value = "'5.041"
$ws.Range("E4").Value = "'-0.50%"
$ws.Range("D5").Value = "'0.08035"
$ws.Range("E5").Value = "'-1.06%"
$ws.Range("D6").Value = "'1.858"
$ws.Range("E6").Value = "'-4.03%"
$ws.Range("D7").Value = "'4.118"
$ws.Range("E7").Value = "'-1.03%"
$ws.Range("D8").Value = "'7.773"
$ws.Range("E8").Value = "'-0.72%"
$ws.Range("D9").Value = "'0.9243"
$ws.Range("E9").Value = "'-1.50%"
$ws.Range("D10").Value = "'0.1268"
$ws.Range("E10").Value = "'-7.23%"
$ws.Range("D11").Value = "'0.1899"
$ws.Range("E11").Value = "'-0.85%"
$ws.Range("D12").Value = "'0.09055"
$ws.Range("E12").Value = "'-1.82%"
$ws.Range("D13").Value = "'0.03435"
$ws.Range("E13").Value = "'-2.03%"
$ws.Range("D14").Value = "'0.09860"
$ws.Range("E14").Value = "'-0.45%"
$ws.Range("D15").Value = "'0.001404"
$ws.Range("E15").Value = "'-1.24%"
$ws.Range("D16").Value = "'0.006258"
$ws.Range("E16").Value = "'7.48%"
$ws.Range("D17").Value = "'3.867"
$ws.Range("E17").Value = "'6.65%"
$ws.Range("E18").Value = "'12.50%"
$ws.Range("D19").Value = "'0.3413"
$ws.Range("E19").Value = "'-0.45%"
$ws.Range("D20").Value = "'0.1334"
$ws.Range("E20").Value = "'-0.82%"
$ws.Range("D21").Value = "'4.796"
$ws.Range("E22").Value = "'-7.71%"
$ws.Range("D23").Value = "'0.04373"
$ws.Range("E23").Value = "'-0.67%"
$ws.Range("E24").Value = "'-0.62%"
$ws.Range("D25").Value = "'0.004857"
$ws.Range("E25").Value = "'1.79%"
$ws.Range("D27").Value = "'0.0001298"
$ws.Range("E27").Value = "'-0.24%"
$ws.Range("E28").Value = "'42.17%"
$ws.Range("D39").Value = "'0.01962"
$ws.Range("E39").Value = "'-2.93%"
$ws.Range("D40").Value = "'0.05157"
$ws.Range("E41").Value = "'-1.30%"
$ws.Range("D42").Value = "'0.01012"
$ws.Range("E42").Value = "'-9.96%"
$ws.Range("D43").Value = "'0.1352"
$ws.Range("E43").Value = "'-2.04%"
$ws.Range("D44").Value = "'0.002107"
$ws.Range("E44").Value = "'0.23%"
$ws.Range("D45").Value = "'0.009876"
$ws.Range("E45").Value = "'-12.58%"
$ws.Range("D46").Value = "'0.00006186"
$ws.Range("E46").Value = "'-2.27%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("D48").Value = "'64.96"
$ws.Range("E48").Value = "'-0.40%"
$ws.Range("D49").Value = "'0.001248"
$ws.Range("E49").Value = "'4.85%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.15%"
